{"js": "// Revature final project doc \u2014 clarification edits:\n//  1) \"last 4 char of ssn +  \"  ->  \"last 6 char of ssn\"\n//     (bump 4 -> 6, and drop the dangling \" +  \" after \"ssn\")\n//  2) \"Calculate the premium end date, premium amount accordingly.\"\n//     -> \"...accordingly. Take the tenure divide into 12 for months.\"\n\nconst body = context.document.body;\n\n// --- Change 1: \"4 char of\" -> \"6 char of\" -------------------------------\nconst charCountHits = body.search(\"4 char of\", { matchCase: true, matchWholeWord: false });\ncharCountHits.load(\"text\");\nawait context.sync();\n\nif (charCountHits.items.length > 0) {\n  charCountHits.items[0].insertText(\"6 char of\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 1b: drop the trailing \" +  \" that used to follow \"ssn\" ------\nconst danglingPlus = body.search(\" +  \", { matchCase: true, matchWholeWord: false });\ndanglingPlus.load(\"text\");\nawait context.sync();\n\nif (danglingPlus.items.length > 0) {\n  danglingPlus.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: extend the \"Calculate the premium...\" sentence -----------\nconst premiumHits = body.search(\"Calculate the premium end date, premium amount accordingly.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\npremiumHits.load(\"text\");\nawait context.sync();\n\nif (premiumHits.items.length > 0) {\n  premiumHits.items[0].insertText(\n    \"Calculate the premium end date, premium amount accordingly. Take the tenure divide into 12 for months.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Revature final project doc - clarification edits:\n#  1) \"last 4 char of ssn +  \"  ->  \"last 6 char of ssn\"\n#     (bump 4 -> 6, and drop the dangling \" +  \" after \"ssn\")\n#  2) \"Calculate the premium end date, premium amount accordingly.\"\n#     -> \"...accordingly. Take the tenure divide into 12 for months.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($FindText, $ReplaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $ReplaceText\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute(\n        [ref]$FindText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $ReplaceText,\n        2\n    )\n}\n\n# --- Change 1: \"4 char of\" -> \"6 char of\" --------------------------------\nReplace-AllText \"4 char of\" \"6 char of\"\n\n# --- Change 1b: drop the trailing \" +  \" that used to follow \"ssn\" -------\nReplace-AllText \" +  \" \"\"\n\n# --- Change 2: extend the \"Calculate the premium...\" sentence -----------\nReplace-AllText \"Calculate the premium end date, premium amount accordingly.\" \"Calculate the premium end date, premium amount accordingly. Take the tenure divide into 12 for months.\"\n"}
